$wb = $excel.ActiveWorkbook

# A cell that is not touched by this edit but already carries the
# "Hyperlink" cell style (s="1") - used below to restore that exact
# style after Hyperlinks.Add() resets it to a freshly generated xf.
$wsStyleSrc = $wb.Worksheets.Item("Register")
$styleSrcCell = $wsStyleSrc.Range("A4")

# ---------------------------------------------------------------
# Sheet "Valid_Login": insert a new leading "scenario" column
# ---------------------------------------------------------------
$wsValid = $wb.Worksheets.Item("Valid_Login")
$wsValid.Range("A1").EntireColumn.Insert()
$wsValid.Range("A1").Value = "scenario"
$wsValid.Range("A2").Value = "validLoginId"
$wsValid.Columns.Item(1).ColumnWidth = 9.89

# Hyperlink used to live on B2, now lives on C2 after the column insert.
$wsValid.Range("B2").Hyperlinks.Delete()
$wsValid.Hyperlinks.Add($wsValid.Range("C2"), "mailto:SeleniumProj@25")
# Re-apply the original (non-hyperlink-default) style to C2.
$styleSrcCell.Copy()
$wsValid.Range("C2").PasteSpecial(-4122)
$wsValid.Range("C2").Value = "SeleniumProj@25"

$wsValid.Range("B1").Select()

# ---------------------------------------------------------------
# Sheet "Login": insert a new leading "scenario" column
# ---------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("A1").EntireColumn.Insert()
$wsLogin.Range("A1").Value = "scenario"
# Fill in "invalidLogin" before "EmptyField" so the shared-string table
# ends up in the same order as the target workbook.
$wsLogin.Range("A4").Value = "invalidLogin"
$wsLogin.Range("A5").Value = "invalidLogin"
$wsLogin.Range("A2").Value = "EmptyField"
$wsLogin.Range("A3").Value = "EmptyField"
$wsLogin.Columns.Item(1).ColumnWidth = 9.8

# Hyperlink used to live on B5, now lives on C5 after the column insert.
$wsLogin.Range("B5").Hyperlinks.Delete()
$wsLogin.Hyperlinks.Add($wsLogin.Range("C5"), "mailto:SeleniumProj@25")
$styleSrcCell.Copy()
$wsLogin.Range("C5").PasteSpecial(-4122)
$wsLogin.Range("C5").Value = "SeleniumProj@25"

$wsLogin.PageSetup.Orientation = 1

$wsLogin.Range("C4").Select()
$wsLogin.Activate()

Write-Host "done"
